$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Scroll the sheet back up so A4 is the top-left visible cell (was showing A22)
[void]$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# Correct the client_id values in column C
$ws.Range("C4").Value = 12345
$ws.Range("C5").Value = 12346
$ws.Range("C6").Value = 12347
$ws.Range("C7").Value = 12345
$ws.Range("C8").Value = 12345
$ws.Range("C9").Value = 12345
$ws.Range("C10").Value = 12347
$ws.Range("C12").Value = 12348
$ws.Range("C13").Value = 12346
$ws.Range("C14").Value = 12347
$ws.Range("C16").Value = 12346
$ws.Range("C18").Value = 12348
$ws.Range("C20").Value = 12347
$ws.Range("C21").Value = 12345
$ws.Range("C22").Value = 12347
$ws.Range("C24").Value = 12347

# Restore the original selection shown in the saved view
[void]$ws.Range("C26").Select()
